$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B4").Value = 7.462000000000001
$ws.Range("B6").Value = 6.034
$ws.Range("B7").Value = 6.032
$ws.Range("D7").Value = -7.647
$ws.Range("D12").Value = -7.302
$ws.Range("E13").Value = 16.468
$ws.Range("E14").Value = 17.176
$ws.Range("D15").Value = -8.411
$ws.Range("B16").Value = 4.935
$ws.Range("E16").Value = 16.409
$ws.Range("E19").Value = 16.54
$ws.Range("B20").Value = 8.466000000000001
$ws.Range("D20").Value = -7.957000000000001
$ws.Range("D21").Value = -8.18
$ws.Range("D22").Value = -7.609999999999999
$ws.Range("E22").Value = 16.627
$ws.Range("D23").Value = -7.87
$ws.Range("B28").Value = 6.171
$ws.Range("B29").Value = 5.101
$ws.Range("D29").Value = -6.927
$ws.Range("B32").Value = 6.572
$ws.Range("D34").Value = -7.917999999999999
$ws.Range("E36").Value = 16.691
$ws.Range("B40").Value = 9.370000000000001
$ws.Range("D42").Value = -8.18
$ws.Range("D43").Value = -7.773000000000001
$ws.Range("D44").Value = -7.860000000000001
$ws.Range("D45").Value = -7.531000000000001
$ws.Range("B46").Value = 5.801
$ws.Range("D46").Value = -7.94
$ws.Range("E46").Value = 16.819
$ws.Range("D50").Value = -8.135
$ws.Range("E50").Value = 16.605
$ws.Range("B51").Value = 4.862
$ws.Range("D51").Value = -8.337
$ws.Range("B52").Value = 5.369
$ws.Range("B57").Value = 5.188999999999999
$ws.Range("B59").Value = 4.714
$ws.Range("B62").Value = 5.229000000000001
$ws.Range("B66").Value = 5.773
$ws.Range("D66").Value = -7.465999999999999
$ws.Range("D67").Value = -7.145
$ws.Range("B73").Value = 6.503
$ws.Range("B74").Value = 9.204000000000001
$ws.Range("D79").Value = -7.642
$ws.Range("D84").Value = -8.300000000000001
$ws.Range("B92").Value = 5.513
$ws.Range("D92").Value = -6.528
$ws.Range("E95").Value = 17.483
$ws.Range("D97").Value = -8.102
$ws.Range("E97").Value = 16.845
$ws.Range("B100").Value = 5.849
